# Sprint Burndown Chart update
# - Sheet "Vokabeltrainer Sprint 1": D4 gets an explicit value of 0
#   (the "Storypoints Plan / done" cell for the first day, which was
#   previously blank).
# - The active selection on that sheet moves from B7 to D5.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Vokabeltrainer Sprint 1")

$ws1.Range("D4").Value = 0

$ws1.Activate() | Out-Null
$ws1.Range("D5").Select() | Out-Null
